# Section 4.3: Degree of lexical flexibility
# Add a new day's worth of writing-progress data (2021-08-31, serial 44212)
# to the "writing" table. Excel auto-grows the table / dimension / named
# range when a ListRow is appended, and the two calculated columns
# (Total / Daily) need their formulas (re)applied for the new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("writing")

$lo = $ws.ListObjects.Item(1)
$newRow = $lo.ListRows.Add()
$r = $newRow.Range

# Raw, directly-entered columns (Date, CV, Abstract ... App: Annotations)
$r.Item(1, 1).Value  = 44212   # A Date
$r.Item(1, 2).Value  = 631     # B CV
$r.Item(1, 3).Value  = 260     # C Abstract
$r.Item(1, 4).Value  = 87      # D Acknowledgments
$r.Item(1, 5).Value  = 229     # E Languages
$r.Item(1, 6).Value  = 492     # F Abbreviations
$r.Item(1, 7).Value  = 416     # G Conventions
$r.Item(1, 8).Value  = 8725    # H 1: Introduction
$r.Item(1, 9).Value  = 19876   # I 2: Background
$r.Item(1, 10).Value = 6562    # J 3: Data & Methods
$r.Item(1, 11).Value = 2870    # K 4: Results
$r.Item(1, 12).Value = 186     # L 5: Conclusion
$r.Item(1, 13).Value = 496     # M App: Grammars
$r.Item(1, 14).Value = 3575    # N App: 100-items
$r.Item(1, 15).Value = 1329    # O App: Annotations

# A54 keeps the same short-date display as the rest of column A
$r.Item(1, 1).NumberFormat = "m/d/yyyy"

# Calculated columns: Total (P) and Daily word-count delta (Q)
$r.Item(1, 16).Formula = "=SUM(C54:O54)"

$qCell = $ws.Range("Q54")
$qCell.Formula = "=SUM(ABS(C54:O54-C53:O53))"

$excel.Calculate()

# Update the dashboard chart's series so they reference the grown range
$dash = $wb.Worksheets.Item("dashboard")
$co = $dash.ChartObjects().Item(1)
$chart = $co.Chart

$s1 = $chart.SeriesCollection().Item(1)
$s1.Formula = "=SERIES(writing!`$Q`$1,writing!`$A`$2:`$A`$54,writing!`$Q`$2:`$Q`$54,1)"

$s2 = $chart.SeriesCollection().Item(2)
$s2.Formula = "=SERIES(writing!`$P`$1,writing!`$A`$2:`$A`$54,writing!`$P`$2:`$P`$54,2)"

# Leave the selection on the newly-entered row (matches the authored edit)
$ws.Range("L54").Select() | Out-Null

# Restore "dashboard" as the active tab (it was active before this edit)
$dash.Activate()
